$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 6 new rows (4-9) continuing the pattern established in rows 2-3.
for ($i = 3; $i -le 8; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = "level 1 6RES source plate"
    $ws.Cells.Item($row, 3).Value = "6RES_AQ_BP"
    $ws.Cells.Item($row, 4).Value = "A1"
    $ws.Cells.Item($row, 5).Value = "384-Well Level 1 MoClo output plate"
    $ws.Cells.Item($row, 6).Value = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"
    $ws.Cells.Item($row, 7).Value = "A$i"
    $ws.Cells.Item($row, 8).Value = 1875
    $ws.Cells.Item($row, 9).Value = "Deionised water"
}
